# tareas_importar.xlsx - "Add files via upload" edit
#
# Adds two new task rows (37 & 38) to Sheet1, widens column B to fit the
# new (longer) text, and moves the selection/scroll position back to the
# top of the sheet (A27) instead of the bottom (A28/A37) it had before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 37: "Terminada" / Reunión dotación Bíobío
# ---------------------------------------------------------------------
$ws.Cells.Item(37, 1).Value = "Terminada"

$ws.Cells.Item(37, 2).Value = "Reunión dotación Bíobío"
$ws.Cells.Item(37, 2).Style = "Normal"

$ws.Cells.Item(37, 3).Value = "Reunión para definir algunos cupos de prioridad"
$ws.Cells.Item(37, 3).Style = "Normal"

$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(37, 5).Value = 1

$ws.Cells.Item(34, 6).Copy()
$ws.Cells.Item(37, 6).PasteSpecial(-4122)
$ws.Cells.Item(37, 6).Value = 45888

$ws.Cells.Item(34, 7).Copy()
$ws.Cells.Item(37, 7).PasteSpecial(-4122)
$ws.Cells.Item(37, 7).Value = 45888

$ws.Cells.Item(34, 8).Copy()
$ws.Cells.Item(37, 8).PasteSpecial(-4122)
$ws.Cells.Item(37, 8).Value = 45888

# ---------------------------------------------------------------------
# Row 38: "Pendiente" / Situación plataforma Linares con cambio de dependencia
# ---------------------------------------------------------------------
$ws.Cells.Item(38, 1).Value = "Pendiente"
$ws.Cells.Item(38, 1).Style = "Normal"

$ws.Cells.Item(38, 2).Value = "Situación plataforma Linares con cambio de dependencia"
$ws.Cells.Item(38, 2).Style = "Normal"

$ws.Cells.Item(38, 3).Value = "Ver los acuerdos con Bastián "
$ws.Cells.Item(38, 3).Style = "Normal"

$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 1

$ws.Cells.Item(34, 6).Copy()
$ws.Cells.Item(38, 6).PasteSpecial(-4122)
$ws.Cells.Item(38, 6).Value = 45888

$ws.Cells.Item(34, 7).Copy()
$ws.Cells.Item(38, 7).PasteSpecial(-4122)
$ws.Cells.Item(38, 7).Value = 45888

# ---------------------------------------------------------------------
# Column B needs to widen to fit the longest new entry (best-fit autosize
# landed at 65.43 chars wide in the real workbook).
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 64.65

# ---------------------------------------------------------------------
# Scroll back to the top of the list and select A27 (instead of the
# bottom-of-sheet view/selection left over from entering the new rows).
# ---------------------------------------------------------------------
$ws.Range("A27").Select()
